# Scheduled-runner data sync: refresh computed profit figures on the
# per-job "Sheets" (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with latest
# market-board pricing. Only literal numeric cells are touched; no
# formulas/formatting are modified.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 8372
$ws.Range("I6").Value = 10026.4
$ws.Range("K6").Value = 30079.2
$ws.Range("M6").Value = -29967.2
$ws.Range("H15").Value = 4193.5815
$ws.Range("I15").Value = 4193.5815
$ws.Range("K15").Value = 12580.7445
$ws.Range("M15").Value = -12411.7445
$ws.Range("H69").Value = 3372.6
$ws.Range("I69").Value = 3245.3333
$ws.Range("J69").Value = 3427.1428
$ws.Range("K69").Value = 9735.999899999999
$ws.Range("L69").Value = 10281.4284
$ws.Range("M69").Value = -8861.999899999999
$ws.Range("N69").Value = -12029.4284
$ws.Range("H72").Value = 3372.6
$ws.Range("I72").Value = 3245.3333
$ws.Range("J72").Value = 3427.1428
$ws.Range("K72").Value = 29207.9997
$ws.Range("L72").Value = 30844.2852
$ws.Range("M72").Value = -24839.9997
$ws.Range("N72").Value = -39580.2852
$ws.Range("H92").Value = 576.3125
$ws.Range("I92").Value = 650.2308
$ws.Range("J92").Value = 256
$ws.Range("K92").Value = 650.2308
$ws.Range("L92").Value = 256
$ws.Range("M92").Value = 597.7692
$ws.Range("N92").Value = -2752
$ws.Range("H112").Value = 2057.4595
$ws.Range("J112").Value = 2084.0557
$ws.Range("L112").Value = 6252.1671
$ws.Range("N112").Value = -8468.167099999999
$ws.Range("H132").Value = 6949889.5
$ws.Range("I132").Value = 10418800
$ws.Range("J132").Value = 12068.75
$ws.Range("K132").Value = 31256400
$ws.Range("L132").Value = 36206.25
$ws.Range("M132").Value = -31253870
$ws.Range("N132").Value = -41266.25
$ws.Range("H137").Value = 1242.3529
$ws.Range("I137").Value = 1201.2142
$ws.Range("J137").Value = 1434.3334
$ws.Range("K137").Value = 3603.6426
$ws.Range("L137").Value = 4303.0002
$ws.Range("M137").Value = -1053.6426
$ws.Range("N137").Value = -9403.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 8000
$ws.Range("J100").Value = 8000
$ws.Range("L100").Value = 8000
$ws.Range("N100").Value = -10164
$ws.Range("H110").Value = 18224.75
$ws.Range("J110").Value = 18224.75
$ws.Range("L110").Value = 18224.75
$ws.Range("N110").Value = -26404.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 939.8570999999999
$ws.Range("I31").Value = 939.8570999999999
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 939.8570999999999
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -644.8570999999999
$ws.Range("H34").Value = 939.8570999999999
$ws.Range("I34").Value = 939.8570999999999
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 939.8570999999999
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -737.8570999999999
$ws.Range("H86").Value = 4779268.5
$ws.Range("J86").Value = 28224.375
$ws.Range("L86").Value = 28224.375
$ws.Range("N86").Value = -30470.375
$ws.Range("H89").Value = 4779268.5
$ws.Range("J89").Value = 28224.375
$ws.Range("L89").Value = 141121.875
$ws.Range("N89").Value = -152353.875
$ws.Range("H92").Value = 32499.334
$ws.Range("J92").Value = 32499.334
$ws.Range("L92").Value = 32499.334
$ws.Range("N92").Value = -37491.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 31.75
$ws.Range("I12").Value = 37.833332
$ws.Range("J12").Value = 25.666666
$ws.Range("K12").Value = 113.499996
$ws.Range("L12").Value = 76.99999800000001
$ws.Range("M12").Value = 59.500004
$ws.Range("N12").Value = -422.999998
$ws.Range("H94").Value = 3837.5
$ws.Range("J94").Value = 4533.3335
$ws.Range("L94").Value = 13600.0005
$ws.Range("N94").Value = -14952.0005
$ws.Range("H122").Value = 776.0833
$ws.Range("I122").Value = 683.2857
$ws.Range("J122").Value = 906
$ws.Range("K122").Value = 6149.571300000001
$ws.Range("L122").Value = 8154
$ws.Range("M122").Value = -3699.571300000001
$ws.Range("N122").Value = -13054
$ws.Range("H126").Value = 5099.968
$ws.Range("J126").Value = 5603.8076
$ws.Range("L126").Value = 16811.4228
$ws.Range("N126").Value = -26691.4228
$ws.Range("H131").Value = 16668006
$ws.Range("J131").Value = 1459.585
$ws.Range("L131").Value = 4378.755
$ws.Range("N131").Value = -14458.755
$ws.Range("H132").Value = 708.8570999999999
$ws.Range("I132").Value = 708.8570999999999
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6379.7139
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -3849.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2339.6428
$ws.Range("I113").Value = 1730.5
$ws.Range("K113").Value = 1730.5
$ws.Range("M113").Value = 439.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 890.65216
$ws.Range("I16").Value = 865.9524
$ws.Range("K16").Value = 865.9524
$ws.Range("M16").Value = -695.9524
$ws.Range("H20").Value = 10000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10452
$ws.Range("H40").Value = 3420.1
$ws.Range("I40").Value = 3139.2
$ws.Range("K40").Value = 3139.2
$ws.Range("M40").Value = -3003.2
$ws.Range("H105").Value = 20000
$ws.Range("J105").Value = 20000
$ws.Range("L105").Value = 20000
$ws.Range("N105").Value = -26988
$ws.Range("H132").Value = 40649.652
$ws.Range("I132").Value = 2127.5715
$ws.Range("J132").Value = 85592.086
$ws.Range("K132").Value = 6382.7145
$ws.Range("L132").Value = 256776.258
$ws.Range("M132").Value = -3852.7145
$ws.Range("N132").Value = -261836.258
$ws.Range("H133").Value = 46442
$ws.Range("J133").Value = 46442
$ws.Range("L133").Value = 46442
$ws.Range("N133").Value = -51502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 2285.7144
$ws.Range("I15").Value = 1500
$ws.Range("K15").Value = 1500
$ws.Range("M15").Value = -1212
$ws.Range("H16").Value = 27000
$ws.Range("J16").Value = 27000
$ws.Range("L16").Value = 27000
$ws.Range("N16").Value = -27584
$ws.Range("H81").Value = 4471.8096
$ws.Range("J81").Value = 5657.3125
$ws.Range("L81").Value = 11314.625
$ws.Range("N81").Value = -13436.625
$ws.Range("H84").Value = 4471.8096
$ws.Range("J84").Value = 5657.3125
$ws.Range("L84").Value = 56573.125
$ws.Range("N84").Value = -67181.125
$ws.Range("H107").Value = 734.3333
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 601.5
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 1804.5
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -5644.5
$ws.Range("H113").Value = 387.62964
$ws.Range("J113").Value = 518.9
$ws.Range("L113").Value = 1556.7
$ws.Range("N113").Value = -5896.7
$ws.Range("H133").Value = 28897.5
$ws.Range("J133").Value = 28897.5
$ws.Range("L133").Value = 28897.5
$ws.Range("N133").Value = -39017.5
